# Update "Horarios" workbook for Línea 141 (LP1912) run at 03:49:28
# - Refresh "Última actualización" timestamp on all three sheets
# - Refresh "Total filas" count and row data on the LP1912 sheet
# - Append three new schedule rows and adjust existing ones

$wb = $excel.ActiveWorkbook

$oldTimestamp = "03:18:49"
$newTimestamp = "03:49:28"

# --- Sheet 1: LP1912 (main schedule data) ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTimestamp"
$ws1.Range("A3").Value = "Total filas: 8"

# Row 6: 14_ABASTO
$ws1.Range("A6").Value = $newTimestamp
$ws1.Range("B6").Value = "03:49"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = "LP1912"

# Row 7: 81_EL PELIGRO
$ws1.Range("A7").Value = $newTimestamp
$ws1.Range("B7").Value = "04:02"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 13
$ws1.Range("E7").Value = "LP1912"

# Row 8: 81_EL PELIGRO
$ws1.Range("A8").Value = $newTimestamp
$ws1.Range("B8").Value = "04:47"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 58
$ws1.Range("E8").Value = "LP1912"

# Row 9: 11_ETCHEVERRY
$ws1.Range("A9").Value = $newTimestamp
$ws1.Range("B9").Value = "04:53"
$ws1.Range("C9").Value = "11_ETCHEVERRY"
$ws1.Range("D9").Value = 64
$ws1.Range("E9").Value = "LP1912"

# Row 10: 17_ROMERO
$ws1.Range("A10").Value = $newTimestamp
$ws1.Range("B10").Value = "05:17"
$ws1.Range("C10").Value = "17_ROMERO"
$ws1.Range("D10").Value = 88
$ws1.Range("E10").Value = "LP1912"

# Row 11 (new): 23_HERNANDEZ
$ws1.Range("A11").Value = $newTimestamp
$ws1.Range("B11").Value = "05:22"
$ws1.Range("C11").Value = "23_HERNANDEZ"
$ws1.Range("D11").Value = 93
$ws1.Range("E11").Value = "LP1912"

# Row 12 (new): 14_ABASTO
$ws1.Range("A12").Value = $newTimestamp
$ws1.Range("B12").Value = "05:43"
$ws1.Range("C12").Value = "14_ABASTO"
$ws1.Range("D12").Value = 114
$ws1.Range("E12").Value = "LP1912"

# Row 13 (new): 17_ROMERO
$ws1.Range("A13").Value = $newTimestamp
$ws1.Range("B13").Value = "05:47"
$ws1.Range("C13").Value = "17_ROMERO"
$ws1.Range("D13").Value = 118
$ws1.Range("E13").Value = "LP1912"

# --- Sheet 2: LP1912-215 (only timestamp refresh) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTimestamp"

# --- Sheet 3: 6203-6173 (only timestamp refresh) ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTimestamp"
